$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.150.24'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '2.363.28'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '547.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.28%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.519'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.75%  '
$ws.Range('D9').Value = '2.368.62'
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.39%  '
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.33'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.350'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.50'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000172'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.35%  '
$ws.Range('D16').Value = '2.798.67'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').Value = '61.340.16'
$ws.Range('E17').Value = '  +1.71%  '
$ws.Range('D18').Value = '2.366.07'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '320.80'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.62'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.29%  '
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.85%  '
$ws.Range('B25').Value = 'SuiNetwork'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.48%  '
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.484.83'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '526.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.34%  '
$ws.Range('D31').Value = '0.0₃0905'
$ws.Range('E31').Value = '  +2.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.39'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.148'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.84'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.50'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.58'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.67'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.89'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.379'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.46'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '146.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.54%  '
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.13%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.74%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.59'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0528'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.580'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0901'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.99%  '
